# Apply the "ig generated files" refresh to the StructureDefinition export:
#   - bump the generation URL from the "pythia" IG to the "cicada" IG
#   - bump the generation Date timestamp
#   - insert a new "Jurisdiction" metadata row (FHIR IGs always emit this
#     row; it was missing before and is empty in this IG)
#   - the newly-inserted row pushes every following metadata row down by
#     one, and the sheet grows from 20 to 21 data rows
#   - the "Elements" sheet references the same generation URL in the
#     Extension.url fixed-value column, so it needs the same text swap

$wb = $excel.ActiveWorkbook

$oldUrl = "http://fhirfli.dev/fhir/ig/pythia/StructureDefinition/allowed-vaccine-status"
$newUrl = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/allowed-vaccine-status"
$newDate = "2026-02-11T14:37:07-05:00"

# ---------------------------------------------------------------------
# "Metadata" sheet: update URL + Date, then make room for the new
# "Jurisdiction" row by shifting rows 11..20 down to rows 12..21.
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = $newUrl
$meta.Range("B8").Value = $newDate

# Shift existing rows 11-20 down into 12-21 (walk bottom-up so we never
# overwrite a row before it has been copied).
for ($r = 20; $r -ge 11; $r--) {
    $destA = $r + 1
    $destB = $r + 1
    $aVal = $meta.Cells.Item($r, 1).Text
    $bVal = $meta.Cells.Item($r, 2).Text
    $meta.Cells.Item($destA, 1).Value = $aVal
    $meta.Cells.Item($destB, 2).Value = $bVal
}

# Row 11 becomes the new "Jurisdiction" property with an empty value.
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# ---------------------------------------------------------------------
# "Elements" sheet: the Extension.url row's Fixed Value column (R5)
# repeats the same generation URL text - keep it in sync.
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("R5").Value = $newUrl
